{"js": "// Replace each math-problem cell's text with its new value.\n// The document body contains a single table laid out as 20 rows x 5\n// columns (100 cells total, row-major). Each cell holds exactly one\n// paragraph with one run whose entire text is a short expression like\n// \"98-69=29\". We walk the table in row-major order and overwrite every\n// cell's paragraph text with the corresponding new expression.\n//\n// We address cells positionally (row, col) instead of searching for the\n// old text, because several of the new expressions equal the *old* text of\n// a different, later cell (e.g. cell index 10 becomes \"52+1=53\", which was\n// originally cell index 74's text) - a sequential search/replace on live\n// text could clobber the wrong cell. Positional addressing sidesteps that.\n//\n// We target each paragraph's own range (not the whole cell body) with\n// insertText(..., \"Replace\") so the existing run formatting (font, size,\n// paragraph alignment) is preserved and only the text itself changes.\nconst newValues = [\n  \"14+2=16\",\n  \"95-77=18\",\n  \"9+18=27\",\n  \"94-72=22\",\n  \"17+27=44\",\n  \"59-23=36\",\n  \"36+14=50\",\n  \"82-35=47\",\n  \"84-54=30\",\n  \"30+25=55\",\n  \"52+1=53\",\n  \"33+47=80\",\n  \"16+55=71\",\n  \"59-9=50\",\n  \"54-27=27\",\n  \"70-45=25\",\n  \"62-0=62\",\n  \"91-62=29\",\n  \"56-49=7\",\n  \"27-13=14\",\n  \"33+48=81\",\n  \"43-16=27\",\n  \"43+7=50\",\n  \"95-59=36\",\n  \"12+65=77\",\n  \"43+33=76\",\n  \"89-59=30\",\n  \"99-36=63\",\n  \"41-10=31\",\n  \"96-87=9\",\n  \"81-26=55\",\n  \"96-31=65\",\n  \"17-14=3\",\n  \"49-30=19\",\n  \"37+24=61\",\n  \"94-87=7\",\n  \"4+84=88\",\n  \"29+50=79\",\n  \"6+16=22\",\n  \"75-5=70\",\n  \"44-40=4\",\n  \"74-49=25\",\n  \"51-38=13\",\n  \"8+54=62\",\n  \"17+5=22\",\n  \"45+1=46\",\n  \"20+43=63\",\n  \"24+47=71\",\n  \"50-16=34\",\n  \"82-73=9\",\n  \"12+64=76\",\n  \"17+61=78\",\n  \"81-10=71\",\n  \"46+26=72\",\n  \"0+16=16\",\n  \"51-44=7\",\n  \"48+1=49\",\n  \"47-13=34\",\n  \"74-49=25\",\n  \"10+6=16\",\n  \"84+0=84\",\n  \"4+71=75\",\n  \"83-41=42\",\n  \"24+37=61\",\n  \"33+45=78\",\n  \"45+34=79\",\n  \"74-15=59\",\n  \"56+38=94\",\n  \"47+46=93\",\n  \"93-75=18\",\n  \"1+78=79\",\n  \"13-9=4\",\n  \"28-0=28\",\n  \"79-47=32\",\n  \"28+63=91\",\n  \"67-18=49\",\n  \"56-27=29\",\n  \"48+29=77\",\n  \"91-67=24\",\n  \"74+14=88\",\n  \"1+59=60\",\n  \"53-43=10\",\n  \"81-0=81\",\n  \"68-15=53\",\n  \"59+40=99\",\n  \"19+67=86\",\n  \"27+69=96\",\n  \"98-97=1\",\n  \"69+1=70\",\n  \"32+16=48\",\n  \"64+13=77\",\n  \"54-6=48\",\n  \"72-37=35\",\n  \"49-7=42\",\n  \"47+14=61\",\n  \"89+6=95\",\n  \"57+14=71\",\n  \"67-64=3\",\n  \"25-3=22\",\n  \"3+16=19\",\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every row's cells up front.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Load every cell's paragraph collection.\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.paragraphs.load(\"items\");\n  }\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const newText = newValues[idx];\n    idx++;\n    const para = cell.body.paragraphs.items[0];\n    const range = para.getRange();\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each math-problem cell's text with its new value.\n# The document body contains a single table laid out as 20 rows x 5\n# columns (100 cells total, row-major). Each cell holds exactly one\n# paragraph with one run whose entire text is a short expression like\n# \"98-69=29\". We walk the table in row-major order (Tables.Item(1).Cell)\n# and overwrite every cell's text with the corresponding new expression.\n#\n# We address cells positionally (row, col) instead of searching for the\n# old text, because several of the new expressions equal the *old* text of\n# a different, later cell (e.g. cell index 10 (1-based 11) becomes\n# \"52+1=53\", which was originally cell index 74's text) - a sequential\n# Find/Replace on live text could clobber the wrong cell. Positional\n# addressing sidesteps that.\n#\n# Setting Cell.Range.Text preserves the existing run/paragraph formatting\n# (font, size, alignment) already present in the cell - only the text\n# itself changes, matching the diff exactly.\n$newValues = @(\n    \"14+2=16\",\n    \"95-77=18\",\n    \"9+18=27\",\n    \"94-72=22\",\n    \"17+27=44\",\n    \"59-23=36\",\n    \"36+14=50\",\n    \"82-35=47\",\n    \"84-54=30\",\n    \"30+25=55\",\n    \"52+1=53\",\n    \"33+47=80\",\n    \"16+55=71\",\n    \"59-9=50\",\n    \"54-27=27\",\n    \"70-45=25\",\n    \"62-0=62\",\n    \"91-62=29\",\n    \"56-49=7\",\n    \"27-13=14\",\n    \"33+48=81\",\n    \"43-16=27\",\n    \"43+7=50\",\n    \"95-59=36\",\n    \"12+65=77\",\n    \"43+33=76\",\n    \"89-59=30\",\n    \"99-36=63\",\n    \"41-10=31\",\n    \"96-87=9\",\n    \"81-26=55\",\n    \"96-31=65\",\n    \"17-14=3\",\n    \"49-30=19\",\n    \"37+24=61\",\n    \"94-87=7\",\n    \"4+84=88\",\n    \"29+50=79\",\n    \"6+16=22\",\n    \"75-5=70\",\n    \"44-40=4\",\n    \"74-49=25\",\n    \"51-38=13\",\n    \"8+54=62\",\n    \"17+5=22\",\n    \"45+1=46\",\n    \"20+43=63\",\n    \"24+47=71\",\n    \"50-16=34\",\n    \"82-73=9\",\n    \"12+64=76\",\n    \"17+61=78\",\n    \"81-10=71\",\n    \"46+26=72\",\n    \"0+16=16\",\n    \"51-44=7\",\n    \"48+1=49\",\n    \"47-13=34\",\n    \"74-49=25\",\n    \"10+6=16\",\n    \"84+0=84\",\n    \"4+71=75\",\n    \"83-41=42\",\n    \"24+37=61\",\n    \"33+45=78\",\n    \"45+34=79\",\n    \"74-15=59\",\n    \"56+38=94\",\n    \"47+46=93\",\n    \"93-75=18\",\n    \"1+78=79\",\n    \"13-9=4\",\n    \"28-0=28\",\n    \"79-47=32\",\n    \"28+63=91\",\n    \"67-18=49\",\n    \"56-27=29\",\n    \"48+29=77\",\n    \"91-67=24\",\n    \"74+14=88\",\n    \"1+59=60\",\n    \"53-43=10\",\n    \"81-0=81\",\n    \"68-15=53\",\n    \"59+40=99\",\n    \"19+67=86\",\n    \"27+69=96\",\n    \"98-97=1\",\n    \"69+1=70\",\n    \"32+16=48\",\n    \"64+13=77\",\n    \"54-6=48\",\n    \"72-37=35\",\n    \"49-7=42\",\n    \"47+14=61\",\n    \"89+6=95\",\n    \"57+14=71\",\n    \"67-64=3\",\n    \"25-3=22\",\n    \"3+16=19\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$numRows = $tbl.Rows.Count\n$numCols = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $numRows; $r++) {\n    for ($c = 1; $c -le $numCols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
